# Update countries & provincias Spain
# - Refresh the "last updated" timestamp.
# - Update case counts for several countries (new data pull).
# - Re-sort a handful of countries whose case totals changed enough to
#   move their rank in the (already totals-sorted) country list:
#     * Panama overtakes Republica Dominicana      (rows 43-44)
#     * Venezuela overtakes Gabon & Guinea          (rows 85-87)
#     * Santa Lucia overtakes Laos                  (rows 203-204, tie - name only)
#     * Groenlandia overtakes Islas Malvinas        (rows 209-210, tie - name only)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados" timestamp (row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 30 de Junio de 2020 a las 03:16"

# --- Estados Unidos (row 4): refreshed totals ---
$ws.Cells.Item(4, 2).Value = 2681802
$ws.Cells.Item(4, 3).Value = 44725
$ws.Cells.Item(4, 4).Value = 1117177
$ws.Cells.Item(4, 5).Value = 1435847
$ws.Cells.Item(4, 7).Value = 341
$ws.Cells.Item(4, 8).Value = 128778

# --- India (row 7): minor recuperados/activos correction ---
$ws.Cells.Item(7, 4).Value = 335271
$ws.Cells.Item(7, 5).Value = 215361

# --- Panama now ranks above Republica Dominicana (rows 43-44) ---
# Row 43 becomes Panama with its freshly updated totals.
$ws.Cells.Item(43, 1).Value = "Panama"
$ws.Cells.Item(43, 2).Value = 32785
$ws.Cells.Item(43, 3).Value = 1099
$ws.Cells.Item(43, 4).Value = 15595
$ws.Cells.Item(43, 5).Value = 16570
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 16
$ws.Cells.Item(43, 8).Value = 620

# Row 44 becomes Republica Dominicana, carrying its prior (unchanged) totals.
$ws.Cells.Item(44, 1).Value = "Republica Dominicana"
$ws.Cells.Item(44, 2).Value = 31816
$ws.Cells.Item(44, 3).Value = 443
$ws.Cells.Item(44, 4).Value = 17280
$ws.Cells.Item(44, 5).Value = 13803
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 7
$ws.Cells.Item(44, 8).Value = 733

# --- Venezuela now ranks above Gabon and Guinea (rows 85-87) ---
# Row 85 becomes Venezuela with its freshly updated totals.
$ws.Cells.Item(85, 1).Value = "Venezuela"
$ws.Cells.Item(85, 2).Value = 5530
$ws.Cells.Item(85, 3).Value = 233
$ws.Cells.Item(85, 4).Value = 1649
$ws.Cells.Item(85, 5).Value = 3833
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 4
$ws.Cells.Item(85, 8).Value = 48

# Row 86 becomes Gabon, carrying its prior (unchanged) totals.
$ws.Cells.Item(86, 1).Value = "Gabon"
$ws.Cells.Item(86, 2).Value = 5394
$ws.Cells.Item(86, 3).Value = 185
$ws.Cells.Item(86, 4).Value = 2420
$ws.Cells.Item(86, 5).Value = 2932
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 2
$ws.Cells.Item(86, 8).Value = 42

# Row 87 becomes Guinea, carrying its prior (unchanged) totals.
$ws.Cells.Item(87, 1).Value = "Guinea"
$ws.Cells.Item(87, 2).Value = 5351
$ws.Cells.Item(87, 3).Value = 9
$ws.Cells.Item(87, 4).Value = 4296
$ws.Cells.Item(87, 5).Value = 1024
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 31

# --- Republica de Africa Central (row 96): refreshed totals ---
$ws.Cells.Item(96, 2).Value = 3613
$ws.Cells.Item(96, 3).Value = 82
$ws.Cells.Item(96, 4).Value = 773
$ws.Cells.Item(96, 5).Value = 2793
$ws.Cells.Item(96, 7).Value = 2
$ws.Cells.Item(96, 8).Value = 47

# --- Guyana (row 166): refreshed totals ---
$ws.Cells.Item(166, 2).Value = 235
$ws.Cells.Item(166, 3).Value = 5
$ws.Cells.Item(166, 4).Value = 114

# --- Bermudas (row 176): refreshed totals ---
$ws.Cells.Item(176, 4).Value = 135
$ws.Cells.Item(176, 5).Value = 2

# --- Santa Lucia now ranks above Laos (rows 203-204, identical totals) ---
$ws.Cells.Item(203, 1).Value = "Santa Lucia"
$ws.Cells.Item(204, 1).Value = "Laos"

# --- Aruba (row 182): refreshed totals ---
$ws.Cells.Item(182, 2).Value = 103
$ws.Cells.Item(182, 3).Value = 2
$ws.Cells.Item(182, 5).Value = 2

# --- Groenlandia now ranks above Islas Malvinas (rows 209-210, identical totals) ---
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"
